# Actualizado taller Robot Delta
# Replace the student roster (column A) with the new list of names and
# update the "grupo" values (column C) that moved between groups.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "BUITRAGO MEDINA, JUAN C.",
    "HERNANDEZ JIMENEZ, ANA M.",
    "JIMENEZ GIRALDO, JUAN P.",
    "MARIN ARROYAVE, JUAN J.",
    "MARIN ZABALA, STEPHANIA",
    "MURIEL AGUDELO, SARA P.",
    "PEÑARANDA BOTELLO, DEIMER L.",
    "ROJAS VASQUEZ, JUAN C.",
    "SOTO JARAMILLO, JUAN J.",
    "TOBON PEÑA, MARIA J.",
    "VALERA MASS, JOSE A.",
    "CASTRO AGUDELO, SAMUEL E.",
    "DE LA CRUZ VERGARA, ELBA LUCIA",
    "GUISAO LOPEZ, JULIAN",
    "HERRERA ESTRADA, SANTIAGO",
    "MEJIA TAMAYO , SANTIAGO",
    "MENESES MONTOYA, SANTIAGO",
    "RAMIREZ VELEZ, CAMILO",
    "TABARES BROWN, TOMAS",
    "VASQUEZ MONTOYA, VALERIA"
)

for ($i = 0; $i -lt $names.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
}

$grupos = @{
    14 = 2
    15 = 2
    18 = 1
    19 = 1
    21 = 1
}

foreach ($row in $grupos.Keys) {
    $ws.Cells.Item($row, 3).Value = $grupos[$row]
}

$ws.Application.GoTo($ws.Range("D16"))
